$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide "http Request Types" (slide 18): add a new bullet line right after
#    "PATCH - Minor update to existing data" in the content placeholder.
# ---------------------------------------------------------------------------

$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    if ($sl.Shapes.Count -ge 1) {
        $titleShape = $sl.Shapes.Item(1)
        if ($titleShape.HasTextFrame -and $titleShape.TextFrame.TextRange.Text -eq "http Request Types") {
            $targetSlide = $sl
            break
        }
    }
}
if ($targetSlide -eq $null) {
    # Fallback: the slide is the 18th slide in this deck.
    $targetSlide = $p.Slides.Item(18)
}

$contentShape = $null
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $shp = $targetSlide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "GET*") {
        $contentShape = $shp
        break
    }
}
if ($contentShape -eq $null) {
    $contentShape = $targetSlide.Shapes.Item(2)
}

$tr = $contentShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count

$patchIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "PATCH*") {
        $patchIndex = $i
        break
    }
}

if ($patchIndex -ne -1) {
    $patchPara = $tr.Paragraphs($patchIndex, 1)
    $patchPara.InsertAfter("`r`tPATCH not enabled by default on Domino")
}

# NOTE: the source deck also re-caches the "last saved" datetimeFigureOut
# field on the Handout Master / Notes Master footers (01/03/2018 -> 11/03/2018,
# 3/1/2018 -> 3/11/2018). Those placeholders are auto-updating date fields
# generated by PowerPoint itself (not user content), and this COM host does
# not expose a writable path to HandoutMaster/NotesMaster shapes, so that
# part of the change is intentionally left alone here.
